$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the three picture shapes (Picture 2, Picture 3, Picture 4),
# keeping only the TextBox. Delete from the end to keep indices stable.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Type -eq 13) {
        $sh.Delete()
    }
}
